$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Week labels in column A: pad single-digit week numbers to two digits (W1 -> W01, ...)
$weekUpdates = @(
    @{ Row = 2;  Week = "W01" },
    @{ Row = 3;  Week = "W02" },
    @{ Row = 4;  Week = "W03" },
    @{ Row = 5;  Week = "W04" },
    @{ Row = 6;  Week = "W05" },
    @{ Row = 7;  Week = "W06" },
    @{ Row = 8;  Week = "W07" },
    @{ Row = 9;  Week = "W08" },
    @{ Row = 10; Week = "W09" }
)

foreach ($u in $weekUpdates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Week
}

# Seasonality Index values in column P
$seasonalityUpdates = @(
    @{ Row = 2;  Value = 1.19 },
    @{ Row = 3;  Value = 1.17 },
    @{ Row = 4;  Value = 0.83 },
    @{ Row = 5;  Value = 1.08 },
    @{ Row = 6;  Value = 1.12 },
    @{ Row = 8;  Value = 0.87 },
    @{ Row = 9;  Value = 0.88 },
    @{ Row = 10; Value = 1.1 },
    @{ Row = 11; Value = 1.18 },
    @{ Row = 12; Value = 1.14 },
    @{ Row = 13; Value = 0.87 },
    @{ Row = 14; Value = 1.03 },
    @{ Row = 15; Value = 0.84 },
    @{ Row = 16; Value = 0.82 },
    @{ Row = 17; Value = 0.89 }
)

foreach ($u in $seasonalityUpdates) {
    $ws.Cells.Item($u.Row, 16).Value = $u.Value
}
